$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'69.363.68"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.72%  '

$ws.Range('D3').Value = "'2.752.20"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +4.25%  '

$ws.Range('D4').Value = "'0.999"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.08%  '

$ws.Range('D5').Value = "'606.05"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.58%  '

$ws.Range('D6').Value = "'169.40"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +6.62%  '

$ws.Range('D7').Value = "'0.999"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.08%  '

$ws.Range('D8').Value = "'0.549"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.06%  '

$ws.Range('D9').Value = "'2.750.70"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.19%  '

$ws.Range('D10').Value = "'0.146"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.36%  '

$ws.Range('D11').Value = "'0.369"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.21%  '

$ws.Range('D12').Value = "'5.38"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.35%  '

$ws.Range('D13').Value = "'0.156"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.08%  '

$ws.Range('D14').Value = "'29.21"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.63%  '

$ws.Range('D15').Value = "'3.252.26"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.25%  '

$ws.Range('D16').Value = "'0.0000192"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.24%  '

$ws.Range('D17').Value = "'69.172.14"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.63%  '

$ws.Range('D18').Value = "'2.756.17"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.86%  '

$ws.Range('D19').Value = "'11.94"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.89%  '

$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').Value = "'372.66"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.50%  '

$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = "'7.77"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +6.12%  '

$ws.Range('D22').Value = "'4.58"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.23%  '

$ws.Range('D23').Value = "'5.04"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.61%  '

$ws.Range('D24').Value = "'2.17"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.55%  '

$ws.Range('D25').Value = "'74.34"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.87%  '

$ws.Range('D27').Value = "'9.93"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.07%  '

$ws.Range('D28').Value = "'2.884.53"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.87%  '

$ws.Range('D29').Value = "'0.0000108"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.76%  '

$ws.Range('D30').Value = "'601.86"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +7.68%  '

$ws.Range('E31').Value = '  +20.47%  '

$ws.Range('D32').Value = "'8.36"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.58%  '

$ws.Range('D33').Value = "'1.47"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.63%  '

$ws.Range('D34').Value = "'1.99"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +7.30%  '

$ws.Range('D35').Value = "'0.135"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.43%  '

$ws.Range('D36').Value = "'1.64"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.82%  '

$ws.Range('D37').Value = "'0.998"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.12%  '

$ws.Range('B38').Value = 'EthereumClassic'
$ws.Range('C38').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D38').Value = "'20.18"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.36%  '

$ws.Range('B39').Value = 'Monero'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D39').Value = "'162.50"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.05%  '

$ws.Range('D40').Value = "'0.386"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.29%  '

$ws.Range('D41').Value = "'1.94"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.99%  '

$ws.Range('D42').Value = "'5.57"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.56%  '

$ws.Range('D43').Value = "'2.76"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.65%  '

$ws.Range('D44').Value = "'18.05"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.45%  '

$ws.Range('D45').Value = "'0.0₆0320"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.46%  '

$ws.Range('E46').Value = '  +0.04%  '

$ws.Range('B47').Value = 'OKB'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D47').Value = "'40.95"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.83%  '

$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').Value = "'158.59"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.74%  '

$ws.Range('B49').Value = 'Filecoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D49').Value = "'3.99"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +6.24%  '

$ws.Range('B50').Value = 'Optimism'
$ws.Range('C50').Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range('D50').Value = "'1.82"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +8.59%  '

$ws.Range('B51').Value = 'ARBITRUM'
$ws.Range('C51').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D51').Value = "'0.614"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +8.78%  '
